{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2025-09-03 Wednesday\", \"2025-09-04 Thursday\"],\n  [\"816\u00f78=102, 0\", \"691\u00f79=76, 7\"],\n  [\"778\u00f79=86, 4\", \"629\u00f75=125, 4\"],\n  [\"935\u00f78=116, 7\", \"552\u00f78=69, 0\"],\n  [\"219\u00f76=36, 3\", \"979\u00f73=326, 1\"],\n  [\"983\u00f76=163, 5\", \"870\u00f73=290, 0\"],\n  [\"563\u00f72=281, 1\", \"947\u00f72=473, 1\"],\n  [\"943\u00f72=471, 1\", \"167\u00f73=55, 2\"],\n  [\"554\u00f75=110, 4\", \"842\u00f72=421, 0\"],\n  [\"356\u00f76=59, 2\", \"674\u00f75=134, 4\"],\n  [\"981\u00f74=245, 1\", \"961\u00f76=160, 1\"],\n  [\"515\u00f78=64, 3\", \"885\u00f74=221, 1\"],\n  [\"159\u00f74=39, 3\", \"402\u00f75=80, 2\"],\n  [\"420\u00f79=46, 6\", \"253\u00f74=63, 1\"],\n  [\"795\u00f77=113, 4\", \"670\u00f77=95, 5\"],\n  [\"640\u00f75=128, 0\", \"359\u00f79=39, 8\"],\n  [\"367\u00f72=183, 1\", \"737\u00f79=81, 8\"],\n  [\"561\u00f78=70, 1\", \"691\u00f75=138, 1\"],\n  [\"554\u00f78=69, 2\", \"297\u00f75=59, 2\"],\n  [\"602\u00f76=100, 2\", \"554\u00f79=61, 5\"],\n  [\"544\u00f79=60, 4\", \"683\u00f73=227, 2\"],\n  [\"250\u00f75=50, 0\", \"604\u00f79=67, 1\"],\n  [\"546\u00f75=109, 1\", \"735\u00f77=105, 0\"],\n  [\"729\u00f74=182, 1\", \"826\u00f76=137, 4\"],\n  [\"744\u00f78=93, 0\", \"273\u00f76=45, 3\"],\n  [\"841\u00f76=140, 1\", \"592\u00f73=197, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load('items');\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const item of found.items) {\n    item.insertText(newText, 'Replace');\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-09-03 Wednesday\", \"2025-09-04 Thursday\")\n    ,@(\"816\u00f78=102, 0\", \"691\u00f79=76, 7\")\n    ,@(\"778\u00f79=86, 4\", \"629\u00f75=125, 4\")\n    ,@(\"935\u00f78=116, 7\", \"552\u00f78=69, 0\")\n    ,@(\"219\u00f76=36, 3\", \"979\u00f73=326, 1\")\n    ,@(\"983\u00f76=163, 5\", \"870\u00f73=290, 0\")\n    ,@(\"563\u00f72=281, 1\", \"947\u00f72=473, 1\")\n    ,@(\"943\u00f72=471, 1\", \"167\u00f73=55, 2\")\n    ,@(\"554\u00f75=110, 4\", \"842\u00f72=421, 0\")\n    ,@(\"356\u00f76=59, 2\", \"674\u00f75=134, 4\")\n    ,@(\"981\u00f74=245, 1\", \"961\u00f76=160, 1\")\n    ,@(\"515\u00f78=64, 3\", \"885\u00f74=221, 1\")\n    ,@(\"159\u00f74=39, 3\", \"402\u00f75=80, 2\")\n    ,@(\"420\u00f79=46, 6\", \"253\u00f74=63, 1\")\n    ,@(\"795\u00f77=113, 4\", \"670\u00f77=95, 5\")\n    ,@(\"640\u00f75=128, 0\", \"359\u00f79=39, 8\")\n    ,@(\"367\u00f72=183, 1\", \"737\u00f79=81, 8\")\n    ,@(\"561\u00f78=70, 1\", \"691\u00f75=138, 1\")\n    ,@(\"554\u00f78=69, 2\", \"297\u00f75=59, 2\")\n    ,@(\"602\u00f76=100, 2\", \"554\u00f79=61, 5\")\n    ,@(\"544\u00f79=60, 4\", \"683\u00f73=227, 2\")\n    ,@(\"250\u00f75=50, 0\", \"604\u00f79=67, 1\")\n    ,@(\"546\u00f75=109, 1\", \"735\u00f77=105, 0\")\n    ,@(\"729\u00f74=182, 1\", \"826\u00f76=137, 4\")\n    ,@(\"744\u00f78=93, 0\", \"273\u00f76=45, 3\")\n    ,@(\"841\u00f76=140, 1\", \"592\u00f73=197, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $r = $d.Content\n    $ok = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
